$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Materialise rows 18-25 (new rows for the extra "Version 2 model"
#    results) by copying an existing, uniformly-styled data row down
#    into them. This gives every new cell the same border/style as
#    the rest of the table before we touch any values.
# ------------------------------------------------------------------
$ws.Range("A3:E3").Copy($ws.Range("A18:E25"))

# ------------------------------------------------------------------
# 2. Undo all existing merges across the whole data block so every
#    cell becomes individually addressable again (and reverts to the
#    plain un-merged style) before we rebuild the new merge layout.
# ------------------------------------------------------------------
$ws.Range("A2:E25").UnMerge()

# ------------------------------------------------------------------
# 3. Write the refreshed values for every row of the table.
# ------------------------------------------------------------------
# Row 2
$ws.Cells.Item(2,1).Value = "dataset_A"
$ws.Cells.Item(2,2).Value = "chan"
$ws.Cells.Item(2,3).Value = "no"
$ws.Cells.Item(2,4).Value = "adapt"
$ws.Cells.Item(2,5).Value = 0.456

# Row 3
$ws.Cells.Item(3,1).Value = ""
$ws.Cells.Item(3,2).Value = ""
$ws.Cells.Item(3,3).Value = ""
$ws.Cells.Item(3,4).Value = "adaptV1"
$ws.Cells.Item(3,5).Value = 0.532

# Row 4
$ws.Cells.Item(4,1).Value = ""
$ws.Cells.Item(4,2).Value = ""
$ws.Cells.Item(4,3).Value = ""
$ws.Cells.Item(4,4).Value = "base"
$ws.Cells.Item(4,5).Value = 0.311

# Row 5
$ws.Cells.Item(5,1).Value = ""
$ws.Cells.Item(5,2).Value = ""
$ws.Cells.Item(5,3).Value = "temp"
$ws.Cells.Item(5,4).Value = "adapt"
$ws.Cells.Item(5,5).Value = 0.425

# Row 6
$ws.Cells.Item(6,1).Value = ""
$ws.Cells.Item(6,2).Value = ""
$ws.Cells.Item(6,3).Value = ""
$ws.Cells.Item(6,4).Value = "adaptV1"
$ws.Cells.Item(6,5).Value = 0.536

# Row 7
$ws.Cells.Item(7,1).Value = ""
$ws.Cells.Item(7,2).Value = ""
$ws.Cells.Item(7,3).Value = ""
$ws.Cells.Item(7,4).Value = "base"
$ws.Cells.Item(7,5).Value = 0.309

# Row 8
$ws.Cells.Item(8,1).Value = ""
$ws.Cells.Item(8,2).Value = "no"
$ws.Cells.Item(8,3).Value = "no"
$ws.Cells.Item(8,4).Value = "adapt"
$ws.Cells.Item(8,5).Value = 0.447

# Row 9
$ws.Cells.Item(9,1).Value = ""
$ws.Cells.Item(9,2).Value = ""
$ws.Cells.Item(9,3).Value = ""
$ws.Cells.Item(9,4).Value = "adaptV1"
$ws.Cells.Item(9,5).Value = 0.306

# Row 10
$ws.Cells.Item(10,1).Value = ""
$ws.Cells.Item(10,2).Value = ""
$ws.Cells.Item(10,3).Value = ""
$ws.Cells.Item(10,4).Value = "base"
$ws.Cells.Item(10,5).Value = 0.253

# Row 11
$ws.Cells.Item(11,1).Value = ""
$ws.Cells.Item(11,2).Value = ""
$ws.Cells.Item(11,3).Value = "temp"
$ws.Cells.Item(11,4).Value = "adapt"
$ws.Cells.Item(11,5).Value = 0.474

# Row 12
$ws.Cells.Item(12,1).Value = ""
$ws.Cells.Item(12,2).Value = ""
$ws.Cells.Item(12,3).Value = ""
$ws.Cells.Item(12,4).Value = "adaptV1"
$ws.Cells.Item(12,5).Value = 0.342

# Row 13
$ws.Cells.Item(13,1).Value = ""
$ws.Cells.Item(13,2).Value = ""
$ws.Cells.Item(13,3).Value = ""
$ws.Cells.Item(13,4).Value = "base"
$ws.Cells.Item(13,5).Value = 0.405

# Row 14
$ws.Cells.Item(14,1).Value = "dataset_B"
$ws.Cells.Item(14,2).Value = "chan"
$ws.Cells.Item(14,3).Value = "no"
$ws.Cells.Item(14,4).Value = "adapt"
$ws.Cells.Item(14,5).Value = 0.525

# Row 15
$ws.Cells.Item(15,1).Value = ""
$ws.Cells.Item(15,2).Value = ""
$ws.Cells.Item(15,3).Value = ""
$ws.Cells.Item(15,4).Value = "adaptV1"
$ws.Cells.Item(15,5).Value = 0.576

# Row 16
$ws.Cells.Item(16,1).Value = ""
$ws.Cells.Item(16,2).Value = ""
$ws.Cells.Item(16,3).Value = ""
$ws.Cells.Item(16,4).Value = "base"
$ws.Cells.Item(16,5).Value = 0.354

# Row 17
$ws.Cells.Item(17,1).Value = ""
$ws.Cells.Item(17,2).Value = ""
$ws.Cells.Item(17,3).Value = "temp"
$ws.Cells.Item(17,4).Value = "adapt"
$ws.Cells.Item(17,5).Value = 0.561

# Row 18
$ws.Cells.Item(18,1).Value = ""
$ws.Cells.Item(18,2).Value = ""
$ws.Cells.Item(18,3).Value = ""
$ws.Cells.Item(18,4).Value = "adaptV1"
$ws.Cells.Item(18,5).Value = 0.593

# Row 19
$ws.Cells.Item(19,1).Value = ""
$ws.Cells.Item(19,2).Value = ""
$ws.Cells.Item(19,3).Value = ""
$ws.Cells.Item(19,4).Value = "base"
$ws.Cells.Item(19,5).Value = 0.534

# Row 20
$ws.Cells.Item(20,1).Value = ""
$ws.Cells.Item(20,2).Value = "no"
$ws.Cells.Item(20,3).Value = "no"
$ws.Cells.Item(20,4).Value = "adapt"
$ws.Cells.Item(20,5).Value = 0.554

# Row 21
$ws.Cells.Item(21,1).Value = ""
$ws.Cells.Item(21,2).Value = ""
$ws.Cells.Item(21,3).Value = ""
$ws.Cells.Item(21,4).Value = "adaptV1"
$ws.Cells.Item(21,5).Value = 0.538

# Row 22
$ws.Cells.Item(22,1).Value = ""
$ws.Cells.Item(22,2).Value = ""
$ws.Cells.Item(22,3).Value = ""
$ws.Cells.Item(22,4).Value = "base"
$ws.Cells.Item(22,5).Value = 0.386

# Row 23
$ws.Cells.Item(23,1).Value = ""
$ws.Cells.Item(23,2).Value = ""
$ws.Cells.Item(23,3).Value = "temp"
$ws.Cells.Item(23,4).Value = "adapt"
$ws.Cells.Item(23,5).Value = 0.558

# Row 24
$ws.Cells.Item(24,1).Value = ""
$ws.Cells.Item(24,2).Value = ""
$ws.Cells.Item(24,3).Value = ""
$ws.Cells.Item(24,4).Value = "adaptV1"
$ws.Cells.Item(24,5).Value = 0.494

# Row 25
$ws.Cells.Item(25,1).Value = ""
$ws.Cells.Item(25,2).Value = ""
$ws.Cells.Item(25,3).Value = ""
$ws.Cells.Item(25,4).Value = "base"
$ws.Cells.Item(25,5).Value = 0.507

# ------------------------------------------------------------------
# 4. Re-establish the merged cells with the new row layout.
# ------------------------------------------------------------------
$ws.Range("A2:A13").Merge()
$ws.Range("A14:A25").Merge()
$ws.Range("B2:B7").Merge()
$ws.Range("B8:B13").Merge()
$ws.Range("B14:B19").Merge()
$ws.Range("B20:B25").Merge()
$ws.Range("C2:C4").Merge()
$ws.Range("C5:C7").Merge()
$ws.Range("C8:C10").Merge()
$ws.Range("C11:C13").Merge()
$ws.Range("C14:C16").Merge()
$ws.Range("C17:C19").Merge()
$ws.Range("C20:C22").Merge()
$ws.Range("C23:C25").Merge()
